$d = $word.ActiveDocument

$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# Append a new blank paragraph right after "Incluyo mi comentario" ...
$blankParaXml = '<w:p xmlns:w="' + $wNs + '"><w:pPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr></w:p>'
$tail = $d.Content
$tail.Collapse(0)
$tail.InsertXML($blankParaXml)

# ... followed by another new paragraph holding the note text.
$noteParaXml = '<w:p xmlns:w="' + $wNs + '"><w:pPr><w:rPr><w:lang w:val="es-ES"/></w:rPr></w:pPr>' +
               '<w:r><w:rPr><w:lang w:val="es-ES"/></w:rPr><w:t>Asigno una nota 8</w:t></w:r></w:p>'
$tail2 = $d.Content
$tail2.Collapse(0)
$tail2.InsertXML($noteParaXml)
